$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price column (D) retains text formatting so values like
# "67.40" or "0.06080" keep trailing zeros instead of becoming numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.468.72'
$ws.Range('E2').Value = '  +1.75%  '
$ws.Range('D3').Value = '1.867.22'
$ws.Range('E3').Value = '  +0.90%  '
$ws.Range('D4').Value = '1.012'
$ws.Range('E4').Value = '  -0.42%  '
$ws.Range('D5').Value = '311.96'
$ws.Range('E5').Value = '  +0.70%  '
$ws.Range('D6').Value = '1.012'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').Value = '0.4784'
$ws.Range('E7').Value = '  +0.26%  '
$ws.Range('D8').Value = '0.3766'
$ws.Range('E8').Value = '  +2.48%  '
$ws.Range('D9').Value = '0.07341'
$ws.Range('E9').Value = '  +1.69%  '
$ws.Range('D10').Value = '0.9362'
$ws.Range('E10').Value = '  +0.80%  '
$ws.Range('D11').Value = '20.71'
$ws.Range('E11').Value = '  +5.07%  '
$ws.Range('D12').Value = '0.07834'
$ws.Range('E12').Value = '  +1.47%  '
$ws.Range('D13').Value = '1.892.62'
$ws.Range('E13').Value = '  +3.45%  '
$ws.Range('D14').Value = '5.438'
$ws.Range('E14').Value = '  +2.13%  '
$ws.Range('D15').Value = '6.557'
$ws.Range('E15').Value = '  +2.05%  '
$ws.Range('D16').Value = '90.51'
$ws.Range('E16').Value = '  +1.84%  '
$ws.Range('D17').Value = '1.013'
$ws.Range('E17').Value = '  -0.40%  '
$ws.Range('D18').Value = '0.000008905'
$ws.Range('E18').Value = '  +3.15%  '
$ws.Range('E19').Value = '  +0.01%  '
$ws.Range('D20').Value = '27.547.09'
$ws.Range('E20').Value = '  +1.97%  '
$ws.Range('D21').Value = '14.75'
$ws.Range('E21').Value = '  +1.55%  '
$ws.Range('D22').Value = '5.118'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').Value = '10.70'
$ws.Range('E23').Value = '  +0.42%  '
$ws.Range('D24').Value = '1.944'
$ws.Range('E24').Value = '  +0.41%  '
$ws.Range('D25').Value = '154.82'
$ws.Range('E25').Value = '  +1.45%  '
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('E27').Value = '  +1.13%  '
$ws.Range('D28').Value = '115.58'
$ws.Range('E28').Value = '  +1.05%  '
$ws.Range('D29').Value = '4.978'
$ws.Range('E29').Value = '  -0.02%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('D31').Value = '3.334'
$ws.Range('E32').Value = '  +3.76%  '
$ws.Range('D33').Value = '0.7564'
$ws.Range('E33').Value = '  +1.86%  '
$ws.Range('D34').Value = '4.604'
$ws.Range('E34').Value = '  +2.33%  '
$ws.Range('D35').Value = '2.738'
$ws.Range('E35').Value = '  +0.63%  '
$ws.Range('D36').Value = '1.121'
$ws.Range('E36').Value = '  +1.07%  '
$ws.Range('D37').Value = '0.02038'
$ws.Range('E37').Value = '  +4.49%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').Value = '0.05262'
$ws.Range('E38').Value = '  -0.14%  '
$ws.Range('B39').Value = 'MXToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D39').Value = '2.990'
$ws.Range('E39').Value = '  +0.30%  '
$ws.Range('D40').Value = '0.5324'
$ws.Range('E40').Value = '  +2.30%  '
$ws.Range('D41').Value = '7.087'
$ws.Range('E41').Value = '  +1.30%  '
$ws.Range('D42').Value = '8.576'
$ws.Range('E42').Value = '  +4.55%  '
$ws.Range('E43').Value = '  +0.95%  '
$ws.Range('D44').Value = '10.66'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('D45').Value = '0.4809'
$ws.Range('E45').Value = '  +1.57%  '
$ws.Range('E46').Value = '  +0.04%  '
$ws.Range('B47').Value = 'NEARProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D47').Value = '1.657'
$ws.Range('E47').Value = '  +3.19%  '
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').Value = '102.87'
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('D49').Value = '67.40'
$ws.Range('E49').Value = '  +2.25%  '
$ws.Range('D50').Value = '0.06080'
$ws.Range('E50').Value = '  +0.79%  '
$ws.Range('D51').Value = '0.9180'
$ws.Range('E51').Value = '  +3.33%  '
